$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.019.55'
$ws.Range('E2').Value = '  +4.10%  '
$ws.Range('D3').Value = '2.650.11'
$ws.Range('E3').Value = '  +6.18%  '
$c = $ws.Range('D4')
$c.Value = "'1.00"
$c.ClearFormats()
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$c.Value = "'114.12"
$c.ClearFormats()
$ws.Range('E5').Value = '  +8.04%  '
$c = $ws.Range('D6')
$c.Value = "'326.75"
$c.ClearFormats()
$ws.Range('E6').Value = '  +2.63%  '
$c = $ws.Range('D7')
$c.Value = "'0.528"
$c.ClearFormats()
$ws.Range('E7').Value = '  +1.81%  '
$ws.Range('E8').Value = '  +0.04%  '
$c = $ws.Range('D9')
$c.Value = "'0.556"
$c.ClearFormats()
$ws.Range('E9').Value = '  +3.77%  '
$c = $ws.Range('D10')
$c.Value = "'40.97"
$c.ClearFormats()
$ws.Range('E10').Value = '  +5.64%  '
$c = $ws.Range('D11')
$c.Value = "'20.20"
$c.ClearFormats()
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('E12').Value = '  +2.67%  '
$ws.Range('E13').Value = '  +0.97%  '
$c = $ws.Range('D14')
$c.Value = "'7.37"
$c.ClearFormats()
$ws.Range('E14').Value = '  +4.26%  '
$ws.Range('D15').Value = '3.065.99'
$ws.Range('E15').Value = '  +6.18%  '
$ws.Range('D16').Value = '2.663.65'
$ws.Range('E16').Value = '  +6.77%  '
$c = $ws.Range('D17')
$c.Value = "'0.875"
$c.ClearFormats()
$ws.Range('E17').Value = '  +5.50%  '
$ws.Range('D18').Value = '49.941.54'
$ws.Range('E18').Value = '  +4.24%  '
$c = $ws.Range('D19')
$c.Value = "'13.23"
$c.ClearFormats()
$ws.Range('E19').Value = '  +2.75%  '
$ws.Range('E20').Value = '  +2.65%  '
$c = $ws.Range('D21')
$c.Value = "'2.92"
$c.ClearFormats()
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').Value = '0.0₃0959'
$ws.Range('E22').Value = '  +3.42%  '
$c = $ws.Range('D23')
$c.Value = "'72.10"
$c.ClearFormats()
$ws.Range('E23').Value = '  +1.54%  '
$c = $ws.Range('D24')
$c.Value = "'277.16"
$c.ClearFormats()
$ws.Range('E24').Value = '  +2.73%  '
$c = $ws.Range('D25')
$c.Value = "'2.59"
$c.ClearFormats()
$ws.Range('E25').Value = '  +3.06%  '
$c = $ws.Range('D26')
$c.Value = "'26.83"
$c.ClearFormats()
$ws.Range('E26').Value = '  +4.06%  '
$ws.Range('E27').Value = '  -0.04%  '
$c = $ws.Range('D28')
$c.Value = "'10.01"
$c.ClearFormats()
$ws.Range('E28').Value = '  +3.18%  '
$ws.Range('E29').Value = '  -2.31%  '
$c = $ws.Range('D30')
$c.Value = "'36.09"
$c.ClearFormats()
$ws.Range('E30').Value = '  +4.68%  '
$ws.Range('E31').Value = '  +2.80%  '
$c = $ws.Range('D32')
$c.Value = "'50.27"
$c.ClearFormats()
$ws.Range('E32').Value = '  +2.03%  '
$ws.Range('E33').Value = '  +3.31%  '
$c = $ws.Range('D34')
$c.Value = "'19.46"
$c.ClearFormats()
$ws.Range('E34').Value = '  +3.08%  '
$c = $ws.Range('D35')
$c.Value = "'0.0810"
$c.ClearFormats()
$ws.Range('E35').Value = '  +4.98%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +7.15%  '
$c = $ws.Range('D38')
$c.Value = "'4.94"
$c.ClearFormats()
$ws.Range('E38').Value = '  +7.98%  '
$ws.Range('E39').Value = '  +8.33%  '
$c = $ws.Range('D40')
$c.Value = "'124.14"
$c.ClearFormats()
$ws.Range('E40').Value = '  +1.91%  '
$ws.Range('E41').Value = '  +2.29%  '
$ws.Range('E42').Value = '  +0.27%  '
$c = $ws.Range('D43')
$c.Value = "'22.17"
$c.ClearFormats()
$ws.Range('E43').Value = '  -0.69%  '
$ws.Range('E44').Value = '  +4.52%  '
$ws.Range('D45').Value = '2.086.96'
$ws.Range('E45').Value = '  +4.45%  '
$c = $ws.Range('D46')
$c.Value = "'3.33"
$c.ClearFormats()
$ws.Range('E46').Value = '  +6.83%  '
$c = $ws.Range('D47')
$c.Value = "'2.33"
$c.ClearFormats()
$ws.Range('E47').Value = '  +16.45%  '
$c = $ws.Range('D48')
$c.Value = "'1.99"
$c.ClearFormats()
$ws.Range('E48').Value = '  +6.25%  '
$c = $ws.Range('D49')
$c.Value = "'9.15"
$c.ClearFormats()
$ws.Range('E49').Value = '  +2.70%  '
$ws.Range('E50').Value = '  +4.89%  '
$c = $ws.Range('D51')
$c.Value = "'59.92"
$c.ClearFormats()
$ws.Range('E51').Value = '  +5.97%  '
